$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Fill in the two previously-empty rows with their new TODO text.
$ws.Range("B28").Value = "validar DocenteCursoABM (validar campos y que el docente no se inscriba 2 veces en un mismo curso)"
$ws.Range("B29").Value = "Armar al menos 2 reportes (se me ocurre ""alumnos con mejor promedio de estado academico"" y ""profesores con mayor carga horaria"")"

# Row 29 grows taller to fit the wrapped, two-line report description (matches ht="30").
$ws.Rows.Item(29).RowHeight = 30

# Update the saved selection to match the author's view (moved on from A27 to B25).
$ws.Range("B25").Select()
